# Edit captured from the authored commit:
#   - the table on slide 5 gets a different built-in table style applied
#     (Table Design gallery -> new tableStyleId GUID)
#   - the deck's theme colour scheme is swapped back to the standard
#     "Office" palette (the file had the "Integral"/"Red Violet" palette
#     applied; PowerPoint re-files the swapped palette across the theme
#     parts on save, so we drive the visible/active colour scheme - the
#     one every slide master/layout/slide actually renders with - back
#     to the plain Office RGB values via the ThemeColorScheme API).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 -------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{2C80D2A6-25DB-496F-9104-B1726D326CFF}")
    }
}

# --- 2. Theme colour scheme -> back to the plain "Office" palette -------
# (COM RGB() packs as R + G*256 + B*65536; values below are precomputed
# from the target #RRGGBB hex so they survive round-tripping untouched.)
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
